# chore: update Sheets via scheduled runner
# Refresh market-board derived price/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) on the per-job Leve tables. Values below were
# recomputed from the latest Universalis price pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1391151.6
$ws.Range("J17").Value = 1483848.4
$ws.Range("L17").Value = 4451545.199999999
$ws.Range("N17").Value = -4451881.199999999

# Row 76
$ws.Range("H76").Value = 13806.308
$ws.Range("J76").Value = 6790
$ws.Range("L76").Value = 6790
$ws.Range("N76").Value = -7420

# Row 79
$ws.Range("H79").Value = 13806.308
$ws.Range("J79").Value = 6790
$ws.Range("L79").Value = 6790
$ws.Range("N79").Value = -8974

# Row 111
$ws.Range("H111").Value = 2252.84
$ws.Range("I111").Value = 2325.3333
$ws.Range("J111").Value = 1872.25
$ws.Range("K111").Value = 6975.999899999999
$ws.Range("L111").Value = 5616.75
$ws.Range("M111").Value = -3908.999899999999
$ws.Range("N111").Value = -11750.75

$ws = $wb.Worksheets.Item("ARM")
# Row 14
$ws.Range("H14").Value = 1269.8
$ws.Range("I14").Value = 1849.6666
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 1849.6666
$ws.Range("L14").Value = 400
$ws.Range("M14").Value = -1674.6666
$ws.Range("N14").Value = -750

# Row 76
$ws.Range("H76").Value = 52666.332
$ws.Range("J76").Value = 52666.332
$ws.Range("L76").Value = 52666.332
$ws.Range("N76").Value = -53342.332

# Row 79
$ws.Range("H79").Value = 52666.332
$ws.Range("J79").Value = 52666.332
$ws.Range("L79").Value = 52666.332
$ws.Range("N79").Value = -55006.332

# Row 80
$ws.Range("H80").Value = 76975
$ws.Range("J80").Value = 76975
$ws.Range("L80").Value = 76975
$ws.Range("N80").Value = -78971

# Row 83
$ws.Range("H83").Value = 76975
$ws.Range("J83").Value = 76975
$ws.Range("L83").Value = 230925
$ws.Range("N83").Value = -240909

# Row 103
$ws.Range("H103").Value = 69000
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 122
$ws.Range("H122").Value = 3277.0833
$ws.Range("I122").Value = 2473
$ws.Range("J122").Value = 4081.1667
$ws.Range("K122").Value = 7419
$ws.Range("L122").Value = 12243.5001
$ws.Range("M122").Value = -4969
$ws.Range("N122").Value = -17143.5001

$ws = $wb.Worksheets.Item("BSM")
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 2650
$ws.Range("J15").Value = 2650
$ws.Range("L15").Value = 2650
$ws.Range("N15").Value = -2990

# Row 44
$ws.Range("H44").Value = 29000
$ws.Range("J44").Value = 29000
$ws.Range("L44").Value = 29000
$ws.Range("N44").Value = -29884

# Row 82
$ws.Range("H82").Value = 49999.5
$ws.Range("J82").Value = 49999.5
$ws.Range("L82").Value = 49999.5
$ws.Range("N82").Value = -50721.5

# Row 85
$ws.Range("H85").Value = 49999.5
$ws.Range("J85").Value = 49999.5
$ws.Range("L85").Value = 49999.5
$ws.Range("N85").Value = -52495.5

# Row 99
$ws.Range("H99").Value = 21824.812
$ws.Range("J99").Value = 7333.6665
$ws.Range("L99").Value = 7333.6665
$ws.Range("N99").Value = -10329.6665

# Row 105
$ws.Range("H105").Value = 8948.286
$ws.Range("I105").Value = 1199.4
$ws.Range("K105").Value = 1199.4
$ws.Range("M105").Value = 547.5999999999999

# Row 126
$ws.Range("H126").Value = 21824.812
$ws.Range("J126").Value = 7333.6665
$ws.Range("L126").Value = 22000.9995
$ws.Range("N126").Value = -26940.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 19
$ws.Range("H19").Value = 2249.8
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 15000
$ws.Range("N19").Value = -15348

# Row 99
$ws.Range("H99").Value = 2228.75
$ws.Range("I99").Value = 1305
$ws.Range("K99").Value = 3915
$ws.Range("M99").Value = -1669

# Row 108
$ws.Range("H108").Value = 2449.5
$ws.Range("I108").Value = 2739.4
$ws.Range("J108").Value = 1000
$ws.Range("K108").Value = 8218.200000000001
$ws.Range("L108").Value = 3000
$ws.Range("M108").Value = -5338.200000000001
$ws.Range("N108").Value = -8760

# Row 109
$ws.Range("H109").Value = 2710.8333
$ws.Range("I109").Value = 2500
$ws.Range("K109").Value = 7500
$ws.Range("M109").Value = -6460

# Row 131
$ws.Range("H131").Value = 1652.6111
$ws.Range("J131").Value = 1641.2572
$ws.Range("L131").Value = 4923.7716
$ws.Range("N131").Value = -15003.7716

# Row 140
$ws.Range("H140").Value = 1174.6
$ws.Range("I140").Value = 968.1111
$ws.Range("K140").Value = 2904.3333
$ws.Range("M140").Value = 2275.6667

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1058
$ws.Range("I22").Value = 375.6
$ws.Range("K22").Value = 375.6
$ws.Range("M22").Value = -80.60000000000002

# Row 27
$ws.Range("H27").Value = 1058
$ws.Range("I27").Value = 375.6
$ws.Range("K27").Value = 375.6
$ws.Range("M27").Value = -268.6

# Row 51
$ws.Range("H51").Value = 39984
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 39984
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 39984
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -40940

# Row 61
$ws.Range("H61").Value = 3138.1904
$ws.Range("I61").Value = 2150.1667
$ws.Range("K61").Value = 2150.1667
$ws.Range("M61").Value = -1948.1667

# Row 81
$ws.Range("H81").Value = 72500
$ws.Range("J81").Value = 64250
$ws.Range("L81").Value = 64250
$ws.Range("N81").Value = -66246

# Row 84
$ws.Range("H84").Value = 72500
$ws.Range("J84").Value = 64250
$ws.Range("L84").Value = 192750
$ws.Range("N84").Value = -202734

# Row 87
$ws.Range("H87").Value = 55499.75
$ws.Range("J87").Value = 56000
$ws.Range("L87").Value = 56000
$ws.Range("N87").Value = -58246

# Row 90
$ws.Range("H90").Value = 55499.75
$ws.Range("J90").Value = 56000
$ws.Range("L90").Value = 168000
$ws.Range("N90").Value = -179232

# Row 113
$ws.Range("H113").Value = 3138.1904
$ws.Range("I113").Value = 2150.1667
$ws.Range("K113").Value = 2150.1667
$ws.Range("M113").Value = 19.83329999999978

# Row 136
$ws.Range("H136").Value = 2073.5
$ws.Range("I136").Value = 1026.3077
$ws.Range("K136").Value = 3078.9231
$ws.Range("M136").Value = -528.9231

$ws = $wb.Worksheets.Item("WVR")
# Row 76
$ws.Range("H76").Value = 46324.332
$ws.Range("J76").Value = 36986.5
$ws.Range("L76").Value = 36986.5
$ws.Range("N76").Value = -37616.5

# Row 79
$ws.Range("H79").Value = 46324.332
$ws.Range("J79").Value = 36986.5
$ws.Range("L79").Value = 36986.5
$ws.Range("N79").Value = -39170.5

# Row 81
$ws.Range("H81").Value = 1521.2
$ws.Range("I81").Value = 1348.3636
$ws.Range("J81").Value = 1996.5
$ws.Range("K81").Value = 2696.7272
$ws.Range("L81").Value = 3993
$ws.Range("M81").Value = -1635.7272
$ws.Range("N81").Value = -6115

# Row 84
$ws.Range("H84").Value = 1521.2
$ws.Range("I84").Value = 1348.3636
$ws.Range("J84").Value = 1996.5
$ws.Range("K84").Value = 13483.636
$ws.Range("L84").Value = 19965
$ws.Range("M84").Value = -8179.635999999999
$ws.Range("N84").Value = -30573

# Row 101
$ws.Range("H101").Value = 14375
$ws.Range("J101").Value = 15833.333
$ws.Range("L101").Value = 15833.333
$ws.Range("N101").Value = -22323.333

# Row 107
$ws.Range("H107").Value = 335.125
$ws.Range("I107").Value = 341.57144
$ws.Range("J107").Value = 290
$ws.Range("K107").Value = 1024.71432
$ws.Range("L107").Value = 870
$ws.Range("M107").Value = 895.28568
$ws.Range("N107").Value = -4710

# Row 113
$ws.Range("H113").Value = 624.2105
$ws.Range("I113").Value = 348.81818
$ws.Range("K113").Value = 1046.45454
$ws.Range("M113").Value = 1123.54546
